# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets,
# reflecting a refreshed data pull (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row => New value for column F, as it appears on sheet "展览" (rows 5-38)
$updatesExhibition = @{
    5  = 15539
    7  = 8
    8  = 701
    9  = 15393
    10 = 51
    11 = 8987
    12 = 375
    14 = 1011
    15 = 88
    18 = 197
    21 = 547
    23 = 10
    24 = 60
    25 = 1109
    26 = 14
    27 = 22
    28 = 80
    30 = 40
    35 = 314
    36 = 449
    38 = 5508
}

# Row => New value for column F, as it appears on sheet "全部类型" (rows 5-40)
$updatesAllTypes = @{
    5  = 15539
    7  = 8
    8  = 701
    9  = 15393
    10 = 51
    11 = 8987
    12 = 375
    14 = 1011
    15 = 88
    18 = 197
    21 = 547
    23 = 10
    24 = 60
    25 = 1109
    26 = 14
    27 = 22
    28 = 80
    30 = 40
    37 = 314
    38 = 449
    40 = 5508
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $updatesExhibition.Keys) {
    $wsExhibition.Range("F$row").Value = $updatesExhibition[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesAllTypes.Keys) {
    $wsAllTypes.Range("F$row").Value = $updatesAllTypes[$row]
}
